$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q3" sheet by copying "2022-Q2" so header/index
#        styling (bold + border) carries over exactly, then overwrite its data. ---
$srcSheet = $wb.Worksheets.Item("2022-Q2")
$srcSheet.Copy($srcSheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# "2022-Q2" has 15 funds (rows 2-16); "2022-Q3" only has 12 (rows 2-13).
# Drop the now-unused trailing rows so the sheet dimension matches.
$newSheet.Range("A14:H16").Clear()

# --- 2. Overwrite the data rows with the 2022-Q3 figures ---
# row 2
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).NumberFormat = "@"
$newSheet.Cells.Item(2,2).Value = "000029"
$newSheet.Cells.Item(2,3).Value = "富国宏观策略灵活配置混合A"
$newSheet.Cells.Item(2,4).NumberFormat = "@"
$newSheet.Cells.Item(2,4).Value = "6.03"
$newSheet.Cells.Item(2,5).NumberFormat = "@"
$newSheet.Cells.Item(2,5).Value = "62.81"
$newSheet.Cells.Item(2,6).NumberFormat = "@"
$newSheet.Cells.Item(2,6).Value = "1.92"
$newSheet.Cells.Item(2,7).NumberFormat = "@"
$newSheet.Cells.Item(2,7).Value = "0.1158"
$newSheet.Cells.Item(2,8).Value = 8

# row 3
$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).NumberFormat = "@"
$newSheet.Cells.Item(3,2).Value = "000849"
$newSheet.Cells.Item(3,3).Value = "汇丰晋信双核策略混合A"
$newSheet.Cells.Item(3,4).NumberFormat = "@"
$newSheet.Cells.Item(3,4).Value = "2.11"
$newSheet.Cells.Item(3,5).NumberFormat = "@"
$newSheet.Cells.Item(3,5).Value = "64.43"
$newSheet.Cells.Item(3,6).NumberFormat = "@"
$newSheet.Cells.Item(3,6).Value = "4.14"
$newSheet.Cells.Item(3,7).NumberFormat = "@"
$newSheet.Cells.Item(3,7).Value = "0.0874"
$newSheet.Cells.Item(3,8).Value = 4

# row 4
$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).NumberFormat = "@"
$newSheet.Cells.Item(4,2).Value = "014175"
$newSheet.Cells.Item(4,3).Value = "工银价值成长混合A"
$newSheet.Cells.Item(4,4).NumberFormat = "@"
$newSheet.Cells.Item(4,4).Value = "2.23"
$newSheet.Cells.Item(4,5).NumberFormat = "@"
$newSheet.Cells.Item(4,5).Value = "66.94"
$newSheet.Cells.Item(4,6).NumberFormat = "@"
$newSheet.Cells.Item(4,6).Value = "2.63"
$newSheet.Cells.Item(4,7).NumberFormat = "@"
$newSheet.Cells.Item(4,7).Value = "0.0586"
$newSheet.Cells.Item(4,8).Value = 8

# row 5
$newSheet.Cells.Item(5,1).Value = 3
$newSheet.Cells.Item(5,2).NumberFormat = "@"
$newSheet.Cells.Item(5,2).Value = "519618"
$newSheet.Cells.Item(5,3).Value = "银河君信灵活配置混合I"
$newSheet.Cells.Item(5,4).NumberFormat = "@"
$newSheet.Cells.Item(5,4).Value = "3.02"
$newSheet.Cells.Item(5,5).NumberFormat = "@"
$newSheet.Cells.Item(5,5).Value = "21.05"
$newSheet.Cells.Item(5,6).NumberFormat = "@"
$newSheet.Cells.Item(5,6).Value = "0.88"
$newSheet.Cells.Item(5,7).NumberFormat = "@"
$newSheet.Cells.Item(5,7).Value = "0.0266"
$newSheet.Cells.Item(5,8).Value = 8

# row 6
$newSheet.Cells.Item(6,1).Value = 4
$newSheet.Cells.Item(6,2).NumberFormat = "@"
$newSheet.Cells.Item(6,2).Value = "014541"
$newSheet.Cells.Item(6,3).Value = "华安新能源主题混合A"
$newSheet.Cells.Item(6,4).NumberFormat = "@"
$newSheet.Cells.Item(6,4).Value = "1.16"
$newSheet.Cells.Item(6,5).NumberFormat = "@"
$newSheet.Cells.Item(6,5).Value = "84.17"
$newSheet.Cells.Item(6,6).NumberFormat = "@"
$newSheet.Cells.Item(6,6).Value = "2.10"
$newSheet.Cells.Item(6,7).NumberFormat = "@"
$newSheet.Cells.Item(6,7).Value = "0.0244"
$newSheet.Cells.Item(6,8).Value = 9

# row 7
$newSheet.Cells.Item(7,1).Value = 5
$newSheet.Cells.Item(7,2).NumberFormat = "@"
$newSheet.Cells.Item(7,2).Value = "005357"
$newSheet.Cells.Item(7,3).Value = "富国国企改革灵活配置混合"
$newSheet.Cells.Item(7,4).NumberFormat = "@"
$newSheet.Cells.Item(7,4).Value = "1.05"
$newSheet.Cells.Item(7,5).NumberFormat = "@"
$newSheet.Cells.Item(7,5).Value = "68.43"
$newSheet.Cells.Item(7,6).NumberFormat = "@"
$newSheet.Cells.Item(7,6).Value = "2.29"
$newSheet.Cells.Item(7,7).NumberFormat = "@"
$newSheet.Cells.Item(7,7).Value = "0.0240"
$newSheet.Cells.Item(7,8).Value = 7

# row 8
$newSheet.Cells.Item(8,1).Value = 6
$newSheet.Cells.Item(8,2).NumberFormat = "@"
$newSheet.Cells.Item(8,2).Value = "000850"
$newSheet.Cells.Item(8,3).Value = "汇丰晋信双核策略混合C"
$newSheet.Cells.Item(8,4).NumberFormat = "@"
$newSheet.Cells.Item(8,4).Value = "0.33"
$newSheet.Cells.Item(8,5).NumberFormat = "@"
$newSheet.Cells.Item(8,5).Value = "64.43"
$newSheet.Cells.Item(8,6).NumberFormat = "@"
$newSheet.Cells.Item(8,6).Value = "4.14"
$newSheet.Cells.Item(8,7).NumberFormat = "@"
$newSheet.Cells.Item(8,7).Value = "0.0137"
$newSheet.Cells.Item(8,8).Value = 4

# row 9
$newSheet.Cells.Item(9,1).Value = 7
$newSheet.Cells.Item(9,2).NumberFormat = "@"
$newSheet.Cells.Item(9,2).Value = "014176"
$newSheet.Cells.Item(9,3).Value = "工银价值成长混合C"
$newSheet.Cells.Item(9,4).NumberFormat = "@"
$newSheet.Cells.Item(9,4).Value = "0.39"
$newSheet.Cells.Item(9,5).NumberFormat = "@"
$newSheet.Cells.Item(9,5).Value = "66.94"
$newSheet.Cells.Item(9,6).NumberFormat = "@"
$newSheet.Cells.Item(9,6).Value = "2.63"
$newSheet.Cells.Item(9,7).NumberFormat = "@"
$newSheet.Cells.Item(9,7).Value = "0.0103"
$newSheet.Cells.Item(9,8).Value = 8

# row 10
$newSheet.Cells.Item(10,1).Value = 8
$newSheet.Cells.Item(10,2).NumberFormat = "@"
$newSheet.Cells.Item(10,2).Value = "013025"
$newSheet.Cells.Item(10,3).Value = "富国宏观策略灵活配置混合C"
$newSheet.Cells.Item(10,4).NumberFormat = "@"
$newSheet.Cells.Item(10,4).Value = "0.35"
$newSheet.Cells.Item(10,5).NumberFormat = "@"
$newSheet.Cells.Item(10,5).Value = "62.81"
$newSheet.Cells.Item(10,6).NumberFormat = "@"
$newSheet.Cells.Item(10,6).Value = "1.92"
$newSheet.Cells.Item(10,7).NumberFormat = "@"
$newSheet.Cells.Item(10,7).Value = "0.0067"
$newSheet.Cells.Item(10,8).Value = 8

# row 11
$newSheet.Cells.Item(11,1).Value = 9
$newSheet.Cells.Item(11,2).NumberFormat = "@"
$newSheet.Cells.Item(11,2).Value = "519617"
$newSheet.Cells.Item(11,3).Value = "银河君信灵活配置混合C"
$newSheet.Cells.Item(11,4).NumberFormat = "@"
$newSheet.Cells.Item(11,4).Value = "0.46"
$newSheet.Cells.Item(11,5).NumberFormat = "@"
$newSheet.Cells.Item(11,5).Value = "21.05"
$newSheet.Cells.Item(11,6).NumberFormat = "@"
$newSheet.Cells.Item(11,6).Value = "0.88"
$newSheet.Cells.Item(11,7).NumberFormat = "@"
$newSheet.Cells.Item(11,7).Value = "0.0040"
$newSheet.Cells.Item(11,8).Value = 8

# row 12
$newSheet.Cells.Item(12,1).Value = 10
$newSheet.Cells.Item(12,2).NumberFormat = "@"
$newSheet.Cells.Item(12,2).Value = "519616"
$newSheet.Cells.Item(12,3).Value = "银河君信灵活配置混合A"
$newSheet.Cells.Item(12,4).NumberFormat = "@"
$newSheet.Cells.Item(12,4).Value = "0.36"
$newSheet.Cells.Item(12,5).NumberFormat = "@"
$newSheet.Cells.Item(12,5).Value = "21.05"
$newSheet.Cells.Item(12,6).NumberFormat = "@"
$newSheet.Cells.Item(12,6).Value = "0.88"
$newSheet.Cells.Item(12,7).NumberFormat = "@"
$newSheet.Cells.Item(12,7).Value = "0.0032"
$newSheet.Cells.Item(12,8).Value = 8

# row 13
$newSheet.Cells.Item(13,1).Value = 11
$newSheet.Cells.Item(13,2).NumberFormat = "@"
$newSheet.Cells.Item(13,2).Value = "014542"
$newSheet.Cells.Item(13,3).Value = "华安新能源主题混合C"
$newSheet.Cells.Item(13,4).NumberFormat = "@"
$newSheet.Cells.Item(13,4).Value = "0.09"
$newSheet.Cells.Item(13,5).NumberFormat = "@"
$newSheet.Cells.Item(13,5).Value = "84.17"
$newSheet.Cells.Item(13,6).NumberFormat = "@"
$newSheet.Cells.Item(13,6).Value = "2.10"
$newSheet.Cells.Item(13,7).NumberFormat = "@"
$newSheet.Cells.Item(13,7).Value = "0.0019"
$newSheet.Cells.Item(13,8).Value = 9

# --- 3. Update the 总计 (summary) sheet: insert the new 2022-Q3 row on top
#        and shift the existing quarters down by one row. ---
$summary = $wb.Worksheets.Item("总计")

# Row 6 is brand new -- copy A5s index-column style (bold + border) onto it first.
$summary.Range("A5").Copy()
$summary.Range("A6").PasteSpecial(-4122)
$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 15
$summary.Range("D6").Value = 4.87

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 19
$summary.Range("D5").Value = 2.69

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 9
$summary.Range("D4").Value = 1.41

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 15
$summary.Range("D3").Value = 0.6899999999999999

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 12
$summary.Range("D2").Value = 0.38

